$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.663.81"
$ws.Range("E2").Value = "  -4.73%  "
$ws.Range("D3").Value = "3.259.19"
$ws.Range("E3").Value = "  -8.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "589.44"
$c.ClearFormats()
$ws.Range("E5").Value = "  -4.50%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "154.73"
$c.ClearFormats()
$ws.Range("E6").Value = "  -11.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.249.80"
$ws.Range("E8").Value = "  -8.19%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.547"
$c.ClearFormats()
$ws.Range("E9").Value = "  -10.80%  "
$ws.Range("E10").Value = "  -12.71%  "
$ws.Range("E11").Value = "  -4.23%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.512"
$c.ClearFormats()
$ws.Range("E12").Value = "  -13.04%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "38.96"
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000247"
$c.ClearFormats()
$ws.Range("E14").Value = "  -10.94%  "
$ws.Range("D15").Value = "3.783.65"
$ws.Range("E15").Value = "  -8.09%  "
$ws.Range("D16").Value = "67.677.90"
$ws.Range("E16").Value = "  -4.72%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "547.45"
$c.ClearFormats()
$ws.Range("E17").Value = "  -10.82%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.258.38"
$ws.Range("E18").Value = "  -8.22%  "
$ws.Range("E19").Value = "  -13.75%  "
$ws.Range("E20").Value = "  -5.68%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "15.38"
$c.ClearFormats()
$ws.Range("E21").Value = "  -13.59%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.772"
$c.ClearFormats()
$ws.Range("E22").Value = "  -13.23%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.87"
$c.ClearFormats()
$ws.Range("E23").Value = "  -14.06%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "86.24"
$c.ClearFormats()
$ws.Range("E24").Value = "  -12.50%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "13.67"
$c.ClearFormats()
$ws.Range("E25").Value = "  -13.26%  "
$ws.Range("E26").Value = "  -0.03%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "3.22"
$c.ClearFormats()
$ws.Range("E27").Value = "  -15.10%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.35"
$c.ClearFormats()
$ws.Range("E28").Value = "  -8.74%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "29.76"
$c.ClearFormats()
$ws.Range("E29").Value = "  -12.13%  "
$ws.Range("E30").Value = "  -17.30%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.77"
$c.ClearFormats()
$ws.Range("E31").Value = "  -9.37%  "
$ws.Range("E32").Value = "  -11.18%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "559.26"
$c.ClearFormats()
$ws.Range("E33").Value = "  -11.53%  "
$ws.Range("E34").Value = "  -18.70%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.84"
$c.ClearFormats()
$ws.Range("E35").Value = "  -15.34%  "
$ws.Range("E36").Value = "  +0.23%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0447"
$c.ClearFormats()
$ws.Range("E37").Value = "  -6.89%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "53.96"
$c.ClearFormats()
$ws.Range("E38").Value = "  -5.42%  "
$ws.Range("E39").Value = "  -14.36%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0856"
$c.ClearFormats()
$ws.Range("E40").Value = "  -14.72%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.128"
$c.ClearFormats()
$ws.Range("E41").Value = "  -12.08%  "
$ws.Range("D42").Value = "2.953.86"
$ws.Range("E42").Value = "  -12.35%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.66"
$c.ClearFormats()
$ws.Range("E43").Value = "  -24.24%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.264"
$c.ClearFormats()
$ws.Range("E44").Value = "  -15.95%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0591"
$ws.Range("E45").Value = "  -20.86%  "
$ws.Range("E46").Value = "  -19.98%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "26.54"
$c.ClearFormats()
$ws.Range("E47").Value = "  -17.65%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.ClearFormats()
$ws.Range("E48").Value = "  -16.13%  "
$ws.Range("E49").Value = "  -0.01%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "125.70"
$c.ClearFormats()
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("E51").Value = "  -12.59%  "
